$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.121.33'
$ws.Range("E2").Value = '  -2.59%  '

$ws.Range("D3").Value = '1.870.51'
$ws.Range("E3").Value = '  -1.87%  '

$ws.Range("D4").Value = '''1.002'

$ws.Range("D5").Value = '''307.39'
$ws.Range("E5").Value = '  -1.77%  '

$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").Value = '''0.5057'
$ws.Range("E7").Value = '  +1.03%  '

$ws.Range("E8").Value = '  -1.47%  '

$ws.Range("D9").Value = '''0.07155'
$ws.Range("E9").Value = '  -1.72%  '

$ws.Range("D10").Value = '''0.8898'
$ws.Range("E10").Value = '  -2.19%  '

$ws.Range("D11").Value = '''20.71'
$ws.Range("E11").Value = '  -0.99%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.867.98'
$ws.Range("E12").Value = '  -1.79%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.07567'
$ws.Range("E13").Value = '  -1.01%  '

$ws.Range("D14").Value = '''5.327'
$ws.Range("E14").Value = '  -3.05%  '

$ws.Range("D15").Value = '''89.35'
$ws.Range("E15").Value = '  -2.79%  '

$ws.Range("E16").Value = '  +0.13%  '

$ws.Range("D17").Value = '''0.000008511'
$ws.Range("E17").Value = '  -2.60%  '

$ws.Range("D18").Value = '''14.14'
$ws.Range("E18").Value = '  -3.21%  '

$ws.Range("E19").Value = '  +0.07%  '

$ws.Range("D20").Value = '27.176.19'
$ws.Range("E20").Value = '  -2.51%  '

$ws.Range("D21").Value = '''5.089'
$ws.Range("E21").Value = '  -1.85%  '

$ws.Range("D22").Value = '2.113.17'
$ws.Range("E22").Value = '  -1.99%  '

$ws.Range("E23").Value = '  -1.76%  '

$ws.Range("D24").Value = '''6.484'
$ws.Range("E24").Value = '  -1.49%  '

$ws.Range("D25").Value = '''151.11'
$ws.Range("E25").Value = '  -1.38%  '

$ws.Range("D26").Value = '''1.842'
$ws.Range("E26").Value = '  -1.85%  '

$ws.Range("D27").Value = '''18.01'
$ws.Range("E27").Value = '  -2.09%  '

$ws.Range("D28").Value = '''2.090'
$ws.Range("E28").Value = '  -5.71%  '

$ws.Range("D29").Value = '''112.95'

$ws.Range("D30").Value = '''4.764'
$ws.Range("E30").Value = '  -3.12%  '

$ws.Range("D31").Value = '''4.687'
$ws.Range("E31").Value = '  -1.82%  '

$ws.Range("D32").Value = '''0.08987'
$ws.Range("E32").Value = '  -0.42%  '

$ws.Range("D33").Value = '''0.05141'
$ws.Range("E33").Value = '  -2.53%  '

$ws.Range("D34").Value = '''3.094'
$ws.Range("E34").Value = '  -3.16%  '

$ws.Range("D35").Value = '''0.7447'
$ws.Range("E35").Value = '  -3.58%  '

$ws.Range("D36").Value = '''1.162'

$ws.Range("D37").Value = '''2.549'
$ws.Range("E37").Value = '  +0.70%  '

$ws.Range("D38").Value = '''0.02033'
$ws.Range("E38").Value = '  -2.45%  '

$ws.Range("D39").Value = '''3.042'
$ws.Range("E39").Value = '  +0.66%  '

$ws.Range("D40").Value = '''1.074'
$ws.Range("E40").Value = '  -1.74%  '

$ws.Range("D41").Value = '''0.5360'
$ws.Range("E41").Value = '  -3.51%  '

$ws.Range("D42").Value = '''6.622'
$ws.Range("E42").Value = '  -4.13%  '

$ws.Range("D43").Value = '''115.01'
$ws.Range("E43").Value = '  +2.82%  '

$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("E45").Value = '  -2.28%  '

$ws.Range("D46").Value = '''0.4646'
$ws.Range("E46").Value = '  -3.93%  '

$ws.Range("D47").Value = '''1.000'
$ws.Range("E47").Value = '  +0.07%  '

$ws.Range("D48").Value = '''10.06'
$ws.Range("E48").Value = '  -5.06%  '

$ws.Range("D49").Value = '''1.572'
$ws.Range("E49").Value = '  -3.77%  '

$ws.Range("D50").Value = '''64.73'
$ws.Range("E50").Value = '  -4.25%  '

$ws.Range("D51").Value = '''36.74'
$ws.Range("E51").Value = '  -0.80%  '
